$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")

# Find the last used row on Sheet1 (column A holds the process id)
$lastRow = $sheet1.Cells.Item($sheet1.Rows.Count, 1).End(-4162).Row
$totalProcess = $lastRow - 1

# Sum the waiting_time (column E) and turn_around_time (column G) columns.
# The source sheet stores these as text, so cast each cell to a number.
$totalWaiting = 0.0
$totalTurn = 0.0
for ($r = 2; $r -le $lastRow; $r++) {
    $totalWaiting += [double]$sheet1.Cells.Item($r, 5).Value2
    $totalTurn += [double]$sheet1.Cells.Item($r, 7).Value2
}

$avgWaiting = [math]::Round($totalWaiting / $totalProcess, 2)
$avgTurn = [math]::Round($totalTurn / $totalProcess, 2)

# Add a new worksheet named "Sheet2" right after Sheet1
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Total Process"
$ws2.Range("B1").Value = $totalProcess

$ws2.Range("A2").Value = "AWT (Average Waiting Time)"
$ws2.Range("B2").Value = $avgWaiting

$ws2.Range("A3").Value = "Total Waiting Time"
$ws2.Range("B3").Value = $totalWaiting

$ws2.Range("A4").Value = "ATAT (Average Turn Around Time)"
$ws2.Range("B4").Value = $avgTurn

$ws2.Range("A5").Value = "Total Turn Around Time"
$ws2.Range("B5").Value = $totalTurn
